$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.223927855491638
$ws.Range("B1").Value = 2.171295404434204
$ws.Range("C1").Value = 4.110167026519775
$ws.Range("D1").Value = 3.106759786605835
$ws.Range("E1").Value = 1.085855007171631
